# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D) for the 633c5ece... row
# (row 5) on both the zh-cn and de-de status sheets with freshly generated
# handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-24 06:44:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-24 06:44:41"
